# Update automatic: dades i banners [2026-02-16 07:20]
# Applies updated values scraped from meteo.cat to the daily summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns a text value to a cell without letting Excel
# auto-convert percentage-looking strings ("85%") into numbers.
function Set-TextValue($range, [string]$text) {
    if ($text -match '^-?\d+(\.\d+)?%$') {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.NumberFormat = "General"
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("E2") "2026-02-16 07:18:29"
Set-TextValue $ws.Range("I2") "3.5 mm"
Set-TextValue $ws.Range("E3") "2026-02-16 07:18:31"
Set-TextValue $ws.Range("I3") "1.1 mm"
Set-TextValue $ws.Range("L3") "48.2 km/h - 234º 6:56 TU"
Set-TextValue $ws.Range("M3") "-0.9 °C 6:59 TU"
Set-TextValue $ws.Range("E4") "2026-02-16 07:18:34"
Set-TextValue $ws.Range("J4") "1014.2 hPa"
Set-TextValue $ws.Range("E5") "2026-02-16 07:18:37"
Set-TextValue $ws.Range("I5") "4.1 mm"
Set-TextValue $ws.Range("E6") "2026-02-16 07:18:39"
Set-TextValue $ws.Range("H6") "85%"
Set-TextValue $ws.Range("J6") "1014.3 hPa"
Set-TextValue $ws.Range("O6") "6.8 °C"
Set-TextValue $ws.Range("E7") "2026-02-16 07:18:42"
Set-TextValue $ws.Range("H7") "56%"
Set-TextValue $ws.Range("L7") "38.5 km/h - 350º 6:54 TU"
Set-TextValue $ws.Range("E8") "2026-02-16 07:18:44"
Set-TextValue $ws.Range("H8") "71%"
Set-TextValue $ws.Range("E9") "2026-02-16 07:18:47"
Set-TextValue $ws.Range("E10") "2026-02-16 07:18:49"
Set-TextValue $ws.Range("E11") "2026-02-16 07:18:51"
Set-TextValue $ws.Range("O11") "0.6 °C"
Set-TextValue $ws.Range("E12") "2026-02-16 07:18:54"
Set-TextValue $ws.Range("E13") "2026-02-16 07:18:56"
Set-TextValue $ws.Range("H13") "90%"
Set-TextValue $ws.Range("J13") "1018.7 hPa"
Set-TextValue $ws.Range("K13") "0.0 MJ/m2"
Set-TextValue $ws.Range("O13") "0.5 °C"
Set-TextValue $ws.Range("E14") "2026-02-16 07:18:58"
Set-TextValue $ws.Range("E15") "2026-02-16 07:19:01"
Set-TextValue $ws.Range("H15") "92%"
Set-TextValue $ws.Range("E16") "2026-02-16 07:19:02"
Set-TextValue $ws.Range("H16") "79%"
Set-TextValue $ws.Range("O16") "-0.3 °C"
Set-TextValue $ws.Range("E17") "2026-02-16 07:19:03"
Set-TextValue $ws.Range("K17") "0.0 MJ/m2"
Set-TextValue $ws.Range("O17") "5.5 °C"
Set-TextValue $ws.Range("E18") "2026-02-16 07:19:04"
Set-TextValue $ws.Range("J18") "1014.7 hPa"
Set-TextValue $ws.Range("O18") "3.8 °C"
Set-TextValue $ws.Range("E19") "2026-02-16 07:19:05"
Set-TextValue $ws.Range("N19") "2.2 °C 6:58 TU"
Set-TextValue $ws.Range("O19") "3.0 °C"
Set-TextValue $ws.Range("E20") "2026-02-16 07:19:06"
Set-TextValue $ws.Range("H20") "94%"
Set-TextValue $ws.Range("E21") "2026-02-16 07:19:07"
Set-TextValue $ws.Range("N21") "2.3 °C 6:59 TU"
Set-TextValue $ws.Range("O21") "4.2 °C"
Set-TextValue $ws.Range("E22") "2026-02-16 07:19:08"
Set-TextValue $ws.Range("L22") "50.8 km/h - 336º 6:32 TU"
Set-TextValue $ws.Range("N22") "-6.5 °C 6:38 TU"
Set-TextValue $ws.Range("E23") "2026-02-16 07:19:11"
Set-TextValue $ws.Range("I23") "2.0 mm"
Set-TextValue $ws.Range("O23") "-0.8 °C"
Set-TextValue $ws.Range("E24") "2026-02-16 07:19:13"
Set-TextValue $ws.Range("H24") "73%"
Set-TextValue $ws.Range("J24") "1017.8 hPa"
Set-TextValue $ws.Range("N24") "9.9 °C 6:58 TU"
Set-TextValue $ws.Range("E25") "2026-02-16 07:19:16"
Set-TextValue $ws.Range("H25") "76%"
Set-TextValue $ws.Range("I25") "0.4 mm"
Set-TextValue $ws.Range("E26") "2026-02-16 07:19:18"
Set-TextValue $ws.Range("E27") "2026-02-16 07:19:20"
Set-TextValue $ws.Range("E28") "2026-02-16 07:19:23"
Set-TextValue $ws.Range("H28") "92%"
Set-TextValue $ws.Range("O28") "2.8 °C"
Set-TextValue $ws.Range("E29") "2026-02-16 07:19:25"
Set-TextValue $ws.Range("K29") "0.0 MJ/m2"
Set-TextValue $ws.Range("N29") "3.7 °C 6:35 TU"
Set-TextValue $ws.Range("O29") "4.7 °C"
Set-TextValue $ws.Range("E30") "2026-02-16 07:19:28"
Set-TextValue $ws.Range("L30") "17.3 km/h - 24º 6:35 TU"
Set-TextValue $ws.Range("E31") "2026-02-16 07:19:30"
Set-TextValue $ws.Range("H31") "57%"
Set-TextValue $ws.Range("N31") "12.3 °C 6:57 TU"
Set-TextValue $ws.Range("O31") "13.7 °C"
Set-TextValue $ws.Range("E32") "2026-02-16 07:19:33"
Set-TextValue $ws.Range("H32") "85%"
Set-TextValue $ws.Range("E33") "2026-02-16 07:19:35"
Set-TextValue $ws.Range("E34") "2026-02-16 07:19:37"
Set-TextValue $ws.Range("M34") "4.1 °C 6:36 TU"
Set-TextValue $ws.Range("O34") "3.3 °C"
Set-TextValue $ws.Range("E35") "2026-02-16 07:19:40"
Set-TextValue $ws.Range("M35") "7.4 °C 6:51 TU"
Set-TextValue $ws.Range("E36") "2026-02-16 07:19:42"
Set-TextValue $ws.Range("H36") "93%"
Set-TextValue $ws.Range("J36") "1014.2 hPa"
Set-TextValue $ws.Range("N36") "4.6 °C 6:35 TU"
Set-TextValue $ws.Range("O36") "6.5 °C"
Set-TextValue $ws.Range("E37") "2026-02-16 07:19:45"
Set-TextValue $ws.Range("J37") "1018.0 hPa"
Set-TextValue $ws.Range("E38") "2026-02-16 07:19:47"
Set-TextValue $ws.Range("H38") "93%"
Set-TextValue $ws.Range("M38") "8.1 °C 6:55 TU"
Set-TextValue $ws.Range("O38") "5.6 °C"
Set-TextValue $ws.Range("E39") "2026-02-16 07:19:49"
Set-TextValue $ws.Range("E40") "2026-02-16 07:19:52"
Set-TextValue $ws.Range("N40") "1.2 °C 6:39 TU"
Set-TextValue $ws.Range("O40") "2.6 °C"
Set-TextValue $ws.Range("E41") "2026-02-16 07:19:54"
Set-TextValue $ws.Range("H41") "53%"
Set-TextValue $ws.Range("J41") "1015.8 hPa"
Set-TextValue $ws.Range("E42") "2026-02-16 07:19:57"
Set-TextValue $ws.Range("O42") "6.1 °C"
Set-TextValue $ws.Range("E43") "2026-02-16 07:19:59"
Set-TextValue $ws.Range("O43") "3.0 °C"
Set-TextValue $ws.Range("E44") "2026-02-16 07:20:02"
Set-TextValue $ws.Range("I44") "2.2 mm"
Set-TextValue $ws.Range("E45") "2026-02-16 07:20:04"
Set-TextValue $ws.Range("I45") "2.3 mm"
Set-TextValue $ws.Range("M45") "3.7 °C 6:52 TU"
Set-TextValue $ws.Range("E46") "2026-02-16 07:20:07"
Set-TextValue $ws.Range("J46") "1018.3 hPa"
